$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1) ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1610.531286
$wsSchedule.Range("F2").Value = 26.62915486111111

# --- Sheet "Detailed" (sheet2) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B9").Value = 73.78127000000001
$wsDetailed.Range("B10").Value = 73.1985
$wsDetailed.Range("B11").Value = 56.98
$wsDetailed.Range("C11").Value = "historical"
$wsDetailed.Range("B12").Value = 62.2571
$wsDetailed.Range("C12").Value = "historical"
$wsDetailed.Range("B13").Value = 80.93711999999999
$wsDetailed.Range("B14").Value = 79.88200999999999
$wsDetailed.Range("B15").Value = 63.65493
$wsDetailed.Range("B16").Value = 50.14936
$wsDetailed.Range("B17").Value = 49.96245
$wsDetailed.Range("B18").Value = 44.39991
$wsDetailed.Range("B20").Value = 36.06
$wsDetailed.Range("B21").Value = 0.009379999999999999
$wsDetailed.Range("B22").Value = -0.3133
$wsDetailed.Range("B23").Value = 36.07
$wsDetailed.Range("B24").Value = 0.51
$wsDetailed.Range("B25").Value = 34.01
$wsDetailed.Range("B28").Value = 0
$wsDetailed.Range("B30").Value = 53.10611
$wsDetailed.Range("B31").Value = 56.98
$wsDetailed.Range("B32").Value = 58.81093
$wsDetailed.Range("B33").Value = 52.47475
$wsDetailed.Range("B34").Value = 56.33086
$wsDetailed.Range("B35").Value = 61.21598
$wsDetailed.Range("B36").Value = 64.01355
$wsDetailed.Range("B37").Value = 37.60586
$wsDetailed.Range("B38").Value = 57.00873
$wsDetailed.Range("B39").Value = 80.02
$wsDetailed.Range("B40").Value = 120.01
$wsDetailed.Range("B41").Value = 158.99
$wsDetailed.Range("B42").Value = 159.6199
$wsDetailed.Range("B43").Value = 111.89625
$wsDetailed.Range("B44").Value = 85.95
$wsDetailed.Range("B45").Value = 78.00005
$wsDetailed.Range("B46").Value = 105.79
$wsDetailed.Range("B47").Value = 64.99985
$wsDetailed.Range("B48").Value = 57.73363
$wsDetailed.Range("B49").Value = 59.4713
